$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.06082199999999999
$ws.Range("I2").Value = 0.1716860072883705
$ws.Range("J2").Value = 0.1716860072883705
$ws.Range("O2").Value = 0.03170629435689145
$ws.Range("P2").Value = 0.03170629435689145
$ws.Range("S2").Value = 0.005443527084044486
$ws.Range("T2").Value = 0.005443527084044486

# Row 3
$ws.Range("G3").Value = 0.06082199999999999
$ws.Range("I3").Value = 0.1716860072883705
$ws.Range("J3").Value = 0.1716860072883705
$ws.Range("M3").Value = 7.341370333333333
$ws.Range("N3").Value = 22.024111
$ws.Range("O3").Value = 0.6001826809091065
$ws.Range("P3").Value = 0.6001826809091064
$ws.Range("Q3").Value = 0.4465168264139999
$ws.Range("R3").Value = 4.018651437726
$ws.Range("S3").Value = 0.1030429681289146
$ws.Range("T3").Value = 0.1030429681289146

# Row 4
$ws.Range("G4").Value = 0.06082199999999999
$ws.Range("I4").Value = 0.1716860072883705
$ws.Range("J4").Value = 0.1716860072883705
$ws.Range("M4").Value = 0.008839999999999999
$ws.Range("N4").Value = 0.02652
$ws.Range("O4").Value = 0.0007227008934757686
$ws.Range("P4").Value = 0.0007227008934757686
$ws.Range("Q4").Value = 0.0005376664799999998
$ws.Range("R4").Value = 0.004838998319999999
$ws.Range("S4").Value = 0.0001240776308645927
$ws.Range("T4").Value = 0.0001240776308645927

# Row 5
$ws.Range("G5").Value = 0.06082199999999999
$ws.Range("I5").Value = 0.1716860072883705
$ws.Range("J5").Value = 0.1716860072883705
$ws.Range("M5").Value = 4.493854666666667
$ws.Range("N5").Value = 13.481564
$ws.Range("O5").Value = 0.3673883238405263
$ws.Range("P5").Value = 0.3673883238405263
$ws.Range("Q5").Value = 0.2733252285359999
$ws.Range("R5").Value = 2.459927056824
$ws.Range("S5").Value = 0.06307543444454682
$ws.Range("T5").Value = 0.06307543444454682

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.293441
$ws.Range("H6").Value = 0.880323
$ws.Range("I6").Value = 0.8283139927116295
$ws.Range("J6").Value = 0.8283139927116295
$ws.Range("O6").Value = 0.03170629435689145
$ws.Range("P6").Value = 0.03170629435689145
$ws.Range("Q6").Value = 0.113804636148
$ws.Range("R6").Value = 1.024241725332
$ws.Range("S6").Value = 0.02626276727284697
$ws.Range("T6").Value = 0.02626276727284697

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.293441
$ws.Range("H7").Value = 0.880323
$ws.Range("I7").Value = 0.8283139927116295
$ws.Range("J7").Value = 0.8283139927116295
$ws.Range("M7").Value = 7.341370333333333
$ws.Range("N7").Value = 22.024111
$ws.Range("O7").Value = 0.6001826809091065
$ws.Range("P7").Value = 0.6001826809091064
$ws.Range("Q7").Value = 2.154259051983666
$ws.Range("R7").Value = 19.388331467853
$ws.Range("S7").Value = 0.4971397127801919
$ws.Range("T7").Value = 0.4971397127801919

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.293441
$ws.Range("H8").Value = 0.880323
$ws.Range("I8").Value = 0.8283139927116295
$ws.Range("J8").Value = 0.8283139927116295
$ws.Range("M8").Value = 0.008839999999999999
$ws.Range("N8").Value = 0.02652
$ws.Range("O8").Value = 0.0007227008934757686
$ws.Range("P8").Value = 0.0007227008934757686
$ws.Range("Q8").Value = 0.00259401844
$ws.Range("R8").Value = 0.02334616596
$ws.Range("S8").Value = 0.0005986232626111759
$ws.Range("T8").Value = 0.0005986232626111759

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.293441
$ws.Range("H9").Value = 0.880323
$ws.Range("I9").Value = 0.8283139927116295
$ws.Range("J9").Value = 0.8283139927116295
$ws.Range("M9").Value = 4.493854666666667
$ws.Range("N9").Value = 13.481564
$ws.Range("O9").Value = 0.3673883238405263
$ws.Range("P9").Value = 0.3673883238405263
$ws.Range("Q9").Value = 1.318681207241333
$ws.Range("R9").Value = 11.868130865172
$ws.Range("S9").Value = 0.3043128893959794
$ws.Range("T9").Value = 0.3043128893959794
